$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = 'We report the first case of transmission of Panton-Valentine leukocidin–producing Staphylococcus aureus to a physician during the resuscitation of an infant with fatal pneumonia.
 The physician exhibited numerous furuncles.
 This case highlights the necessity for health care workers to protect themselves against transmission of infectious diseases from patient to care giver.
'
$ws.Range('E2').Value = '[Martin%Chalumeau%NULL%0, Philippe%Bidet%NULL%1, Gérard%Lina%NULL%1, Mostapha%Mokhtari%NULL%1, Marie-Claude%André%NULL%1, Dominique%Gendrel%NULL%1, Edouard%Bingen%NULL%1, Josette%Raymond%j.raymond@svp.ap-hop-paris.fr%1]'
$ws.Range('I2').Value = ''
$ws.Range('J2').Value = 'The University of Chicago Press'

# Row 3
$ws.Range('E3').Value = '[Hae-Sung%Nam%NULL%0, Mi-Yeon%Yeon%NULL%0, Mi-Yeon%Yeon%NULL%0, Jung Wan%Park%NULL%0, Jung Wan%Park%NULL%0, Jee-Young%Hong%NULL%2, Jee-Young%Hong%NULL%0, Ji Woong%Son%NULL%2, Ji Woong%Son%NULL%0]'
$ws.Range('I3').Value = ''
$ws.Range('J3').Value = 'Korean Society of Epidemiology'

# Row 4
$ws.Range('D4').Value = 'Infection of healthcare workers with the severe acute respiratory syndrome–associated coronavirus (SARS-CoV) is thought to occur primarily by either contact or large respiratory droplet transmission.
 However, infrequent healthcare worker infections occurred despite the use of contact and droplet precautions, particularly during certain aerosol-generating medical procedures.
 We investigated a possible cluster of SARS-CoV infections in healthcare workers who used contact and droplet precautions during attempted cardiopulmonary resuscitation of a SARS patient.
 Unlike previously reported instances of transmission during aerosol-generating procedures, the index case-patient was unresponsive, and the intubation procedure was performed quickly and without difficulty.
 However, before intubation, the patient was ventilated with a bag-valve-mask that may have contributed to aerosolization of SARS-CoV.
 On the basis of the results of this investigation and previous reports of SARS transmission during aerosol-generating procedures, a systematic approach to the problem is outlined, including the use of the following: 1) administrative controls, 2) environmental engineering controls, 3) personal protective equipment, and 4) quality control.
'
$ws.Range('E4').Value = '[Michael D.%Christian%NULL%0, Mona%Loutfy%NULL%1, L. Clifford%McDonald%NULL%2, Kenneth F.%Martinez%NULL%1, Mariana%Ofner%NULL%1, Tom%Wong%NULL%1, Tamara%Wallington%NULL%1, Wayne L.%Gold%NULL%1, Barbara%Mederski%NULL%1, Karen%Green%NULL%3, Donald E.%Low%NULL%1, NULL%NULL%NULL%0]'
$ws.Range('I4').Value = ''
$ws.Range('J4').Value = 'Centers for Disease Control and Prevention'

# Row 5
$ws.Range('C5').Value = 'Unknown Title'
$ws.Range('E5').Value = '[]'
$ws.Range('F5').Value = 'not found'
$ws.Range('G5').Value = 'N/A'
$ws.Range('H5').Value = '1970-01-01'
$ws.Range('J5').Value = ''

# Row 6
$ws.Range('C6').Value = 'Unknown Title'
$ws.Range('D6').Value = 'Unknown Abstract'
$ws.Range('E6').Value = '[]'
$ws.Range('F6').Value = 'not found'
$ws.Range('G6').Value = 'N/A'
$ws.Range('H6').Value = '1970-01-01'
$ws.Range('I6').Value = ''

# Row 7
$ws.Range('D7').Value = '
Objective  To evaluate possible severe acute respiratory syndrome (SARS) infection associated risk factors in a SARS affected hospital in Beijing by means of a case control study.
'
$ws.Range('E7').Value = '[Wei%Liu%NULL%0, Fang%Tang%NULL%2, Li‐Qun%Fang%NULL%1, Sake J.%De Vlas%NULL%1, Huai‐Jian%Ma%NULL%1, Jie‐Ping%Zhou%NULL%1, Caspar W. N.%Looman%NULL%1, Jan Hendrik%Richardus%NULL%1, Wu‐Chun%Cao%NULL%1]'
$ws.Range('I7').Value = ''
$ws.Range('J7').Value = 'Blackwell Publishing Ltd'

# Row 8
$ws.Range('D8').Value = 'To determine factors that predispose or protect healthcare workers from severe acute respiratory syndrome (SARS), we conducted a retrospective cohort study among 43 nurses who worked in two Toronto critical care units with SARS patients.
 Eight of 32 nurses who entered a SARS patient’s room were infected.
 The probability of SARS infection was 6% per shift worked.
 Assisting during intubation, suctioning before intubation, and manipulating the oxygen mask were high-risk activities.
 Consistently wearing a mask (either surgical or particulate respirator type N95) while caring for a SARS patient was protective for the nurses, and consistent use of the N95 mask was more protective than not wearing a mask.
 Risk was reduced by consistent use of a surgical mask, but not significantly.
 Risk was lower with consistent use of a N95 mask than with consistent use of a surgical mask.
 We conclude that activities related to intubation increase SARS risk and use of a mask (particularly a N95 mask) is protective.
'
$ws.Range('E8').Value = '[Mark%Loeb%NULL%0, Allison%McGeer%NULL%2, Bonnie%Henry%NULL%2, Marianna%Ofner%NULL%2, David%Rose%NULL%1, Tammy%Hlywka%NULL%1, Joanne%Levie%NULL%1, Jane%McQueen%NULL%1, Stephanie%Smith%NULL%1, Lorraine%Moss%NULL%1, Andrew%Smith%NULL%1, Karen%Green%NULL%0, Stephen D.%Walter%NULL%1]'
$ws.Range('I8').Value = ''
$ws.Range('J8').Value = 'Centers for Disease Control and Prevention'

# Row 9
$ws.Range('E9').Value = '[Janet%Raboud%NULL%0, Altynay%Shigayeva%NULL%1, Allison%McGeer%NULL%0, Erika%Bontovics%NULL%1, Martin%Chapman%NULL%1, Denise%Gravel%NULL%1, Bonnie%Henry%NULL%0, Stephen%Lapinsky%NULL%1, Mark%Loeb%NULL%0, L. Clifford%McDonald%NULL%0, Marianna%Ofner%NULL%0, Shirley%Paton%NULL%1, Donna%Reynolds%NULL%1, Damon%Scales%NULL%1, Sandy%Shen%NULL%1, Andrew%Simor%NULL%1, Thomas%Stewart%NULL%1, Mary%Vearncombe%NULL%1, Dick%Zoutman%NULL%1, Karen%Green%NULL%0, Joel Mark%Montgomery%NULL%4, Joel Mark%Montgomery%NULL%0]'
$ws.Range('I9').Value = ''
$ws.Range('J9').Value = 'Public Library of Science'

# Row 10
$ws.Range('C10').Value = 'Unknown Title'
$ws.Range('D10').Value = 'Unknown Abstract'
$ws.Range('E10').Value = '[]'
$ws.Range('F10').Value = 'not found'
$ws.Range('G10').Value = 'N/A'
$ws.Range('H10').Value = '1970-01-01'
$ws.Range('J10').Value = ''

# Row 11
$ws.Range('E11').Value = '[Hyungoo%Shin%NULL%0, Jaehoon%Oh%NULL%1, Tae Ho%Lim%NULL%1, Hyunggoo%Kang%NULL%1, Yeongtak%Song%NULL%1, Sanghyun%Lee%NULL%1, Abdelouahab%Bellou.%NULL%2, Abdelouahab%Bellou.%NULL%0]'
$ws.Range('I11').Value = ''
$ws.Range('J11').Value = 'Wolters Kluwer Health'

# Row 12
$ws.Range('C12').Value = 'Unknown Title'
$ws.Range('F12').Value = 'not found'
$ws.Range('G12').Value = 'N/A'
$ws.Range('I12').Value = ''

Write-Host "Updated 49 cells"